$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.202.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.105.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.102.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  +6.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.617.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.252.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.112.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "502.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.292.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.32%  "
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0803"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.13%  "
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.256"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.19%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0538"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.59%  "
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("E51").Value = "  +4.21%  "
